# Reproduce the commit "Actualización modelo_af 2022 en ceros":
# Set the full year 2022 daily values (column B, rows 2-366) to 0
# (dates 44562..44926, i.e. 2022-01-01 .. 2022-12-31), and move the
# active selection to D365 as it was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 366 correspond to every day of 2022 (365 rows).
$range = $ws.Range("B2:B366")
$range.Value = 0

# Update the saved selection/active cell, matching the source workbook.
$ws.Range("D365").Select()
